$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.Formula = "'" + $val
    $cell.Style = "Normal"
}

Set-TextCell "D2" '37.443.89'
Set-TextCell "E2" '  -1.23%  '

Set-TextCell "D3" '2.058.68'
Set-TextCell "E3" '  -1.24%  '

Set-TextCell "E4" '  +0.00%  '

Set-TextCell "D5" '231.74'
Set-TextCell "E5" '  -0.72%  '

Set-TextCell "E6" '  -0.59%  '

Set-TextCell "D8" '57.32'
Set-TextCell "E8" '  -3.47%  '

Set-TextCell "D9" '0.386'
Set-TextCell "E9" '  -2.47%  '

Set-TextCell "D10" '0.0774'
Set-TextCell "E10" '  -1.76%  '

Set-TextCell "E11" '  +1.43%  '

Set-TextCell "D12" '2.359.62'
Set-TextCell "E12" '  -1.32%  '

Set-TextCell "D13" '14.66'
Set-TextCell "E13" '  -0.70%  '

Set-TextCell "D14" '21.01'
Set-TextCell "E14" '  -1.06%  '

Set-TextCell "D15" '0.761'
Set-TextCell "E15" '  -2.11%  '

Set-TextCell "D16" '5.34'
Set-TextCell "E16" '  -0.23%  '

Set-TextCell "D17" '2.058.53'
Set-TextCell "E17" '  -1.31%  '

Set-TextCell "D18" '37.379.92'
Set-TextCell "E18" '  -1.13%  '

Set-TextCell "D19" '6.10'
Set-TextCell "E19" '  -0.48%  '

Set-TextCell "D20" '69.76'
Set-TextCell "E20" '  -2.71%  '

Set-TextCell "D21" '0.0₃0826'
Set-TextCell "E21" '  -2.64%  '

Set-TextCell "D22" '226.94'
Set-TextCell "E22" '  -0.52%  '

Set-TextCell "E23" '  +0.07%  '

Set-TextCell "E24" '  +0.33%  '

Set-TextCell "D26" '9.90'
Set-TextCell "E26" '  +7.57%  '

Set-TextCell "D27" '170.51'
Set-TextCell "E27" '  -0.91%  '

Set-TextCell "E28" '  -5.40%  '

Set-TextCell "D29" '19.24'
Set-TextCell "E29" '  -1.36%  '

Set-TextCell "D30" '1.35'
Set-TextCell "E30" '  -4.68%  '

Set-TextCell "E31" '  +0.34%  '

Set-TextCell "E32" '  -3.92%  '

Set-TextCell "D33" '0.0622'
Set-TextCell "E33" '  -1.51%  '

Set-TextCell "D34" '4.61'
Set-TextCell "E34" '  -2.28%  '

Set-TextCell "D35" '2.51'
Set-TextCell "E35" '  -0.01%  '

Set-TextCell "E36" '  +0.40%  '

Set-TextCell "D37" '3.29'
Set-TextCell "E37" '  -3.79%  '

Set-TextCell "E38" '  +0.12%  '

Set-TextCell "E39" '  -1.88%  '

Set-TextCell "D40" '0.0226'
Set-TextCell "E40" '  +3.41%  '

Set-TextCell "D41" '98.45'
Set-TextCell "E41" '  -0.63%  '

Set-TextCell "B42" 'TrustWalletToken'
Set-TextCell "C42" 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell "D42" '1.21'
Set-TextCell "E42" '  +4.50%  '

Set-TextCell "B43" 'Cronos'
Set-TextCell "C43" 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell "D43" '0.0959'
Set-TextCell "E43" '  -2.46%  '

Set-TextCell "D44" '1.477.72'
Set-TextCell "E44" '  +2.19%  '

Set-TextCell "E45" '  -0.52%  '

Set-TextCell "D46" '16.73'
Set-TextCell "E46" '  -0.84%  '

Set-TextCell "E47" '  -2.74%  '

Set-TextCell "B48" 'FTXToken'
Set-TextCell "C48" 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextCell "D48" '3.99'
Set-TextCell "E48" '  -4.66%  '

Set-TextCell "B49" 'FraxShare'
Set-TextCell "C49" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell "D49" '7.25'
Set-TextCell "E49" '  -1.80%  '

Set-TextCell "E50" '  -1.30%  '

Set-TextCell "D51" '2.246.49'
Set-TextCell "E51" '  -1.34%  '
